$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells - copy the formatting (bold, border, centered) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J
$data = @{
    2  = @(8, 9)
    3  = @(6, 7)
    4  = @(8, 8)
    5  = @(9, 9)
    6  = @(7, 8)
    7  = @(7, 8)
    8  = @(4, 5)
    9  = @(7, 8)
    10 = @(6, 7)
    11 = @(7, 8)
    12 = @(8, 8)
    13 = @(6, 7)
    14 = @(1, 1)
    15 = @(1, 5)
    16 = @(1, 7)
    17 = @(1, 7)
    18 = @(1, 7)
    19 = @(1, 5)
    20 = @(1, 5)
    21 = @(1, 5)
    22 = @(1, 7)
    23 = @(1, 5)
    24 = @(1, 3)
    25 = @(6, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
